# Update the weekly report workbook:
# - Refresh the "Report Generated On" timestamp
# - Zero out billed amount / pricing totals (no-violation / no-billable scenario)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Work Report")

# Report generated timestamp
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:49 PM"

# Total Billed Amount
$ws.Range("C8").Value = 0

# Pricing for line item and TOTAL row
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 0
